$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B (old B -> C, old C -> D), keep formatting from the left
$ws.Columns("B").Insert()

# Approximate the original column A width (75.81640625 XML units) for the new column B.
# This runtime's ColumnWidth setter quantizes to 1/6 character-width steps, so we
# compensate by the fixed offset (5/6) this runtime applies when converting the
# "characters" value to the stored XML width.
$ws.Columns("B").ColumnWidth = (75.81640625 - (5/6))

# New header + query text for the stat-bar query column
$ws.Range("B1").Value = "StatQuery"

$ws.Range("B2").Value = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Bernese Mountain Dog']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

# Match wrap-text style of the other long-text cell in this row (A2)
$ws.Range("B2").WrapText = $true

# Move the active selection/view as in the edited workbook
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("B2").Select() | Out-Null
